$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 3): four new columns K, L, M, N
# (written in this exact order so the shared-string table gets indices
#  17="Some dates", 18="Some errors", 19="Dates with NA", 20="Errors with NA")
# ---------------------------------------------------------------------
$ws.Range("K3").Value = "Some dates"
$ws.Range("M3").Value = "Some errors"
$ws.Range("L3").Value = "Dates with NA"
$ws.Range("N3").Value = "Errors with NA"

# ---------------------------------------------------------------------
# Column K: "Some dates"
# ---------------------------------------------------------------------
$ws.Range("K4").Value = 42066
$ws.Range("K4").NumberFormat = "mm-dd-yy"

$ws.Range("K5").Value = 42039.426388888889
$ws.Range("K5").NumberFormat = "m/d/yy h:mm"

$ws.Range("K6").Value = 32242
$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$ws.Range("K7").Value = 0.62638888888888888
$ws.Range("K7").NumberFormat = "[`$-F400]h:mm:ss\ AM/PM"

# ---------------------------------------------------------------------
# Column L: "Dates with NA"
# ---------------------------------------------------------------------
$ws.Range("L4").Value = 23835
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

$ws.Range("L5").Value = 18484.777777777777
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)

$ws.Range("L6").Value = 0.79166666666666663
$ws.Range("L6").NumberFormat = "h:mm"

# ---------------------------------------------------------------------
# Column M: "Some errors"
# ---------------------------------------------------------------------
$ws.Range("M4").Formula = "=1/0"
$ws.Range("M5").Formula = "=NA()"

# Establish the external-workbook link (via a throwaway cell using the
# full '[Book]Sheet'!Ref syntax) so the [1]nonexistingsheet!A1 short form
# resolves to a real externalReference/externalLink part on save.
$ws.Range("ZZ1").Formula = "='[nonexistingfile.xlsx]nonexistingsheet'!A1"
$ws.Range("M6").Formula = "=[1]nonexistingsheet!A1"
$ws.Range("ZZ1").ClearContents()

$ws.Range("M7").Formula = "=a0"

# ---------------------------------------------------------------------
# Column N: "Errors with NA"
# ---------------------------------------------------------------------
$ws.Range("N4").Formula = "=1/0"
$ws.Range("N5").Formula = "=NA()"
$ws.Range("N6").Formula = "=a0"

# ---------------------------------------------------------------------
# Column widths (best-fit) for the whole used range, to match the
# "select columns, AutoFit" step the author performed.
# ---------------------------------------------------------------------
$ws.Range("C:N").EntireColumn.AutoFit()

$excel.CutCopyMode = $false

$ws.Range("N7").Select()
